# TCMD_YR_FIN.xlsx update
# Commit: "Doing Updates for Financials"
# A new most-recent fiscal-year column is inserted as column D (pushing the
# previously-existing D:K data right to E:L) across the three statement
# blocks (Income Statement, Balance Sheet, Cash Flow Statement), and the new
# column is populated with the latest period's figures. One prior data point
# (Capital Expenditures, row 91) is also corrected in the same edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Insert a new blank column at D; existing D:K shift right to E:L.
$ws.Columns("D:D").Insert()

# 2. Pick up the correct number formatting/style for the new column by
#    copying it from the (now-shifted) column E, one sub-table at a time so
#    the stray label-only rows (37, 79) don't pick up a spurious blank cell.
$ws.Range("E7:E35").Copy()
$ws.Range("D7:D35").PasteSpecial(-4122)

$ws.Range("E38:E77").Copy()
$ws.Range("D38:D77").PasteSpecial(-4122)

$ws.Range("E80:E102").Copy()
$ws.Range("D80:D102").PasteSpecial(-4122)

$ws.Application.CutCopyMode = 0

# 3. Populate the new column D with the latest period's values.

# -- Income Statement --
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 143800
$ws.Range("D9").Value = 75100
$ws.Range("D10").Value = 68700
$ws.Range("D12").Value = 5300
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("D17").Value = 140800
$ws.Range("D18").Value = 3000
$ws.Range("D20").Value = 500
$ws.Range("D21").Value = 7100
$ws.Range("D22").Value = 0
$ws.Range("D23").Value = 3500
$ws.Range("D24").Value = -3100
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 6600
$ws.Range("D27").Value = 6600
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = -500
$ws.Range("D33").Value = 6600
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 6600

# -- Balance Sheet --
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 20100
$ws.Range("D42").Value = 25800
$ws.Range("D43").Value = 26100
$ws.Range("D44").Value = 11200
$ws.Range("D45").Value = 1800
$ws.Range("D46").Value = 85000
$ws.Range("D47").Value = 1900
$ws.Range("D48").Value = 4800
$ws.Range("D49").Value = 5300
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 10100
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 107100
$ws.Range("D57").Value = 5100
$ws.Range("D58").Value = "NA"
$ws.Range("D59").Value = 11000
$ws.Range("D60").Value = 16100
$ws.Range("D61").Value = 0
$ws.Range("D62").Value = 1700
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 17800
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 9700
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 89300
$ws.Range("D77").Value = 0

# -- Cash Flow Statement --
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 6600
$ws.Range("D83").Value = 3600
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 9000
$ws.Range("D91").Value = -4200
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -14700
$ws.Range("D96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = 1900
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = -3900

# 4. A prior-period data correction that rode along with this edit: the
#    Capital Expenditures figure that used to sit in (old) D91 moved to E91
#    but with a corrected value (-3700 instead of the old -3800).
$ws.Range("E91").Value = -3700

# 5. Column width / default row height refresh that came with the new column
#    (Excel nudged these slightly when the sheet was resaved).
$ws.Columns("A:A").ColumnWidth = 7.88671875
$ws.Columns("B:B").ColumnWidth = 26.88671875
$ws.Columns("C:C").ColumnWidth = 69.109375
$ws.Range("D1:I1").EntireColumn.ColumnWidth = 14.6640625
$ws.Range("J1:K1").EntireColumn.ColumnWidth = 4.88671875
